# The workbook originally holds a single sheet ("Лист1") that actually
# contains two separate receipt tables stacked in the same sheet: one in
# rows 1-21 and a second one (its own header-less set, numbered 1-20 in
# column A) in rows 24-43, separated by two blank rows (22-23).
#
# The edit splits these into two named sheets:
#   - "Vika"   : keeps the first table (rows 1-21)
#   - "Artyom" : gets the second table, renumbered to rows 1-20
#
# It also updates each sheet's view/selection state to match what Excel
# saved afterwards.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# Grab the second table's values (A24:K43) before touching anything else.
$secondTable = $ws1.Range("A24:K43").Value2

# Rename the existing sheet.
$ws1.Name = "Vika"

# Insert the new sheet right after "Vika" and name it.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Artyom"

# Move the second table onto the new sheet, shifted up to start at row 1.
$ws2.Range("A1:K20").Value2 = $secondTable

# Remove the now-duplicated rows (the 2 blank separator rows plus the
# second table) from the original sheet.
$ws1.Range("A22:K43").ClearContents() | Out-Null

# Match the saved selection/view state: "Vika" ends up the active sheet
# with cell D30 selected (and no leftover scrolled top-left cell), while
# "Artyom" has P19 selected.
$ws2.Range("P19").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D30").Select() | Out-Null
